$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between rows 183 and 184 ---
# (Index column A, and B:E metadata columns, stay untouched for these two rows.)
$row183 = @{}
$row184 = @{}
for ($col = 6; $col -le 22; $col++) {
    $row183[$col] = $ws.Cells.Item(183, $col).Value()
    $row184[$col] = $ws.Cells.Item(184, $col).Value()
}
for ($col = 6; $col -le 22; $col++) {
    $ws.Cells.Item(183, $col).Value = $row184[$col]
    $ws.Cells.Item(184, $col).Value = $row183[$col]
}

# --- Append new row 185 (Reading v Bristol Rovers) ---
# Copy formatting from the row above (184) first, so styles (bold index column,
# date/time number format on column E, etc.) match the rest of the table.
$ws.Range("A184:V184").Copy()
$ws.Range("A185:V185").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(185, 1).Value = 184
$ws.Cells.Item(185, 2).Value = "england"
$ws.Cells.Item(185, 3).Value = "league-one"
$ws.Cells.Item(185, 4).Value = "2023-2024"
$ws.Cells.Item(185, 5).Value = 45237.875
$ws.Cells.Item(185, 6).Value = "Reading"
$ws.Cells.Item(185, 7).Value = 1
$ws.Cells.Item(185, 8).Value = "Bristol Rovers"
$ws.Cells.Item(185, 9).Value = 1
$ws.Cells.Item(185, 10).Value = 3.15
$ws.Cells.Item(185, 11).Value = "30/10/2023 21:42"
$ws.Cells.Item(185, 12).Value = 2.64
$ws.Cells.Item(185, 13).Value = "07/11/2023 20:50"
$ws.Cells.Item(185, 14).Value = 3.49
$ws.Cells.Item(185, 15).Value = "30/10/2023 21:42"
$ws.Cells.Item(185, 16).Value = 3.68
$ws.Cells.Item(185, 17).Value = "07/11/2023 20:50"
$ws.Cells.Item(185, 18).Value = 2.19
$ws.Cells.Item(185, 19).Value = "30/10/2023 21:42"
$ws.Cells.Item(185, 20).Value = 2.61
$ws.Cells.Item(185, 21).Value = "07/11/2023 20:50"
$ws.Cells.Item(185, 22).Value = "https://www.betexplorer.com/football/england/league-one/reading-bristol-rovers/n900clND/"

Write-Host "done"
